# Update the "all electricity sources" logit exponent on the ETLE sheet
# from -5 to -3, and leave that sheet active/selected (matching the
# state the workbook was saved in).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ETLE")

$ws.Range("B2").Value = -3

# Make ETLE the active sheet (mirrors tabSelected/activeTab moving to ETLE)
$ws.Activate()
